$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New data row (row 56) in the FWHM table
$ws.Range("A56").Value = "sg_rr_52_025 2023-12-11 18-30-06.csv"
$ws.Range("B56").Value = 0.01
$ws.Range("C56").Value = 1000
$ws.Range("D56").Value = 5001
$ws.Range("E56").Value = 1530
$ws.Range("F56").Value = 1570
$ws.Range("G56").Value = 0.5
$ws.Range("H56").Value = "(approx_fsr/2)/wavelength step size"
$ws.Range("I56").Value = 2
$ws.Range("J56").Value = 1.8875
$ws.Range("K56").Value = 0.0078430324425366096
$ws.Range("L56").Value = "yes"
$ws.Range("M56").Value = 0.140490273155437
$ws.Range("N56").Value = 0.0066112105211042198
$ws.Range("O56").Value = "reduced approx fsr a bit, to see if this had any affect on fsr calculation as above, half the approx fsr was quite close to actual calculated fsr."

# Update the view so the new row is visible (matches author's saved view state)
$ws.Application.ActiveWindow.ScrollRow = 44
$ws.Application.ActiveWindow.Zoom = 73
$ws.Range("A57").Select()
